$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.131.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07392"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8802"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.28"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.896.10"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07689"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.39"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.375"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.535"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008730"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.578.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.62"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.246"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.875"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.12"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.134"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.179"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.60"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08928"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7440"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.164"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.516"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.940"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.091"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05303"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01936"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.308"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.932"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5262"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1644"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.387"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4907"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.44"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.43"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06278"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.40%  "

# Row 36/37 swap: RenderToken <-> Frax positions
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.008"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.583"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.60%  "
